# feat: photograph file name input feature
#
# The "work experience" block on the "회원 정보" sheet (rows 4-8) gets
# collapsed: the old placeholder row "e" (row 5) is replaced by the real
# column headers (근무기간/근무처/담당업무/근속연수), the old header row
# (row 7) and the trailing placeholder row "t" (row 8) are dropped
# entirely, leaving just the three placeholder/header rows 4-6.
#
# The free-text note on "자기소개서" (sheet 2, A1) is replaced with a
# shorter note and word-wrap is turned on for it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Row 5 ("e" placeholder row) becomes the real header row.
$ws1.Range("A5").Value2 = "근무기간"
$ws1.Range("B5").Value2 = "근무처"
$ws1.Range("C5").Value2 = "담당업무"
$ws1.Range("D5").Value2 = "근속연수"

# Row 6 becomes the "r" placeholder row.
$ws1.Range("A6").Value2 = "r"
$ws1.Range("B6").Value2 = "r"
$ws1.Range("C6").Value2 = "r"
$ws1.Range("D6").Value2 = "r"

# Old rows 7 (근무기간/근무처/담당업무/근속연수 header) and 8 ("t" row)
# are removed completely, shrinking the sheet's used range to A1:F6.
$ws1.Rows("7:8").Delete()

# Sheet2 A1: shorten the free-text note and enable word-wrap on it.
$ws2.Range("A1").Value2 = "asdfasdf`n"
$ws2.Range("A1").WrapText = $true
